# Auto-generated Excel COM-interop script
# Applies timetable updates: adds Lesson 6 / Lesson 7 columns (G,H) and
# replaces lesson text in columns B-H for all 6 sheets (10A,10B,10C,11A,11B,11C)

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---- Sheet '10A' ----
$ws = $wb.Worksheets.Item('10A')

# Add Lesson 6 / Lesson 7 header cells, copying the header style from F1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G1").Value = "Lesson 6"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H1").Value = "Lesson 7"

$values = @{
    'B2' = 'Литература - Room 105 - Teacher: Владимир В.П.'
    'C2' = 'Русский язык - Room 205 - Teacher: Екатерина А.П.'
    'D2' = 'Физика - Room 110 - Teacher: Роман А.Д.'
    'E2' = 'Математика - Room 205 - Teacher: Евгения О.П.'
    'F2' = 'Математика - Room 108 - Teacher: Алексей С.К.'
    'G2' = 'Математика - Room 205 - Teacher: Константин В.Л.'
    'H2' = 'География - Room 109 - Teacher: Марина Г.Б.'
    'B3' = 'Физика - Room 114 - Teacher: Игорь В.Н.'
    'C3' = 'Информатика - Room 107 - Teacher: Анастасия А.М.'
    'D3' = 'Физика - Room 113 - Teacher: Игорь В.Н.'
    'E3' = 'Иностранный язык - Room 207 - Teacher: Елена В.К.'
    'F3' = 'Русский язык - Room 213 - Teacher: Александр И.Р.'
    'G3' = 'Химия - Room 211 - Teacher: Оксана М.С.'
    'H3' = 'Физическая культура - Room 111 - Teacher: Юлия Н.В.'
    'B4' = 'Физическая культура - Room 112 - Teacher: Юлия Н.В.'
    'C4' = 'Обществознание - Room 213 - Teacher: Наталья Д.З.'
    'D4' = 'Русский язык - Room 108 - Teacher: Ольга В.С.'
    'E4' = 'Информатика - Room 107 - Teacher: Игорь В.Н.'
    'F4' = 'Химия - Room 211 - Teacher: Елена С.Т.'
    'G4' = 'Математика - Room 205 - Teacher: Игорь П.В.'
    'H4' = 'История - Room 114 - Teacher: Елена В.К.'
    'B5' = 'Обществознание - Room 214 - Teacher: Наталья Д.З.'
    'C5' = 'Биология - Room 309 - Teacher: Андрей П.С.'
    'D5' = 'История - Room 111 - Teacher: Елена В.К.'
    'E5' = 'Математика - Room 109 - Teacher: Марина Г.Б.'
    'F5' = 'География - Room 109 - Teacher: Василий И.М.'
    'G5' = 'Литература - Room 202 - Teacher: Мария А.К.'
    'H5' = 'Информатика - Room 107 - Teacher: Алёна В.Л.'
    'B6' = 'Русский язык - Room 206 - Teacher: Анна С.И.'
    'C6' = 'Математика - Room 109 - Teacher: Игорь П.В.'
    'D6' = 'Иностранный язык - Room 203 - Teacher: Юлия Н.В.'
    'E6' = 'X - Room 0 - Teacher: X'
    'F6' = 'Литература - Room 207 - Teacher: Светлана А.К.'
    'G6' = 'Химия - Room 212 - Teacher: Анастасия А.Ф.'
    'H6' = 'Биология - Room 309 - Teacher: Светлана А.К.'
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet '10B' ----
$ws = $wb.Worksheets.Item('10B')

# Add Lesson 6 / Lesson 7 header cells, copying the header style from F1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G1").Value = "Lesson 6"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H1").Value = "Lesson 7"

$values = @{
    'B2' = 'Информатика - Room 107 - Teacher: Игорь В.Н.'
    'C2' = 'Русский язык - Room 110 - Teacher: Анна С.И.'
    'D2' = 'Иностранный язык - Room 203 - Teacher: Елена В.К.'
    'E2' = 'Биология - Room 309 - Teacher: Андрей П.С.'
    'F2' = 'История - Room 114 - Teacher: Елена В.К.'
    'G2' = 'X - Room 0 - Teacher: X'
    'H2' = 'Химия - Room 211 - Teacher: Анастасия А.Ф.'
    'B3' = 'Физическая культура - Room 112 - Teacher: Светлана А.К.'
    'C3' = 'Математика - Room 212 - Teacher: Игорь П.В.'
    'D3' = 'Информатика - Room 107 - Teacher: Алёна В.Л.'
    'E3' = 'Русский язык - Room 101 - Teacher: Екатерина А.П.'
    'F3' = 'География - Room 109 - Teacher: Василий И.М.'
    'G3' = 'Русский язык - Room 213 - Teacher: Екатерина А.П.'
    'H3' = 'Математика - Room 104 - Teacher: Марина Г.Б.'
    'B4' = 'Физика - Room 113 - Teacher: Роман А.Д.'
    'C4' = 'Обществознание - Room 212 - Teacher: Наталья Д.З.'
    'D4' = 'История - Room 110 - Teacher: Елена В.К.'
    'E4' = 'Литература - Room 103 - Teacher: Светлана А.К.'
    'F4' = 'Литература - Room 202 - Teacher: Илья В.М.'
    'G4' = 'Математика - Room 205 - Teacher: Константин В.Л.'
    'H4' = 'Физика - Room 110 - Teacher: Людмила А.С.'
    'B5' = 'Математика - Room 212 - Teacher: Сергей А.Т.'
    'C5' = 'Физическая культура - Room 112 - Teacher: Юлия Н.В.'
    'D5' = 'Химия - Room 209 - Teacher: Елена С.Т.'
    'E5' = 'Математика - Room 205 - Teacher: Игорь П.В.'
    'F5' = 'Физика - Room 114 - Teacher: Игорь В.Н.'
    'G5' = 'Обществознание - Room 212 - Teacher: Анна С.И.'
    'H5' = 'Информатика - Room 107 - Teacher: Игорь В.Н.'
    'B6' = 'Русский язык - Room 110 - Teacher: Александр И.Р.'
    'C6' = 'Математика - Room 205 - Teacher: Сергей А.Т.'
    'D6' = 'Литература - Room 209 - Teacher: Георгий Н.М.'
    'E6' = 'Иностранный язык - Room 206 - Teacher: Юлия Н.В.'
    'F6' = 'Химия - Room 212 - Teacher: Оксана М.С.'
    'G6' = 'Биология - Room 310 - Teacher: Андрей П.С.'
    'H6' = 'География - Room 109 - Teacher: Марина Г.Б.'
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet '10C' ----
$ws = $wb.Worksheets.Item('10C')

# Add Lesson 6 / Lesson 7 header cells, copying the header style from F1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G1").Value = "Lesson 6"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H1").Value = "Lesson 7"

$values = @{
    'B2' = 'Физическая культура - Room 112 - Teacher: Юлия Н.В.'
    'C2' = 'Русский язык - Room 206 - Teacher: Александр И.Р.'
    'D2' = 'Математика - Room 212 - Teacher: Сергей А.Т.'
    'E2' = 'Физика - Room 113 - Teacher: Роман А.Д.'
    'F2' = 'Русский язык - Room 213 - Teacher: Михаил С.К.'
    'G2' = 'Математика - Room 205 - Teacher: Игорь П.В.'
    'H2' = 'Иностранный язык - Room 203 - Teacher: Юлия Н.В.'
    'B3' = 'Биология - Room 309 - Teacher: Светлана А.К.'
    'C3' = 'Русский язык - Room 213 - Teacher: Михаил С.К.'
    'D3' = 'Литература - Room 105 - Teacher: Илья В.М.'
    'E3' = 'История - Room 110 - Teacher: Елена В.К.'
    'F3' = 'Информатика - Room 107 - Teacher: Игорь В.Н.'
    'G3' = 'Химия - Room 213 - Teacher: Елена С.Т.'
    'H3' = 'Физика - Room 113 - Teacher: Игорь В.Н.'
    'B4' = 'Математика - Room 109 - Teacher: Константин В.Л.'
    'C4' = 'Обществознание - Room 214 - Teacher: Анна С.И.'
    'D4' = 'Биология - Room 309 - Teacher: Светлана А.К.'
    'E4' = 'Русский язык - Room 206 - Teacher: Ольга В.С.'
    'F4' = 'Иностранный язык - Room 206 - Teacher: Елена В.К.'
    'G4' = 'География - Room 109 - Teacher: Василий И.М.'
    'H4' = 'X - Room 0 - Teacher: X'
    'B5' = 'Физическая культура - Room 111 - Teacher: Юлия Н.В.'
    'C5' = 'Математика - Room 205 - Teacher: Марина Г.Б.'
    'D5' = 'Обществознание - Room 214 - Teacher: Анна С.И.'
    'E5' = 'Математика - Room 108 - Teacher: Евгения О.П.'
    'F5' = 'Литература - Room 105 - Teacher: Петр И.С.'
    'G5' = 'История - Room 111 - Teacher: Наталья Д.З.'
    'H5' = 'Химия - Room 213 - Teacher: Анастасия А.Ф.'
    'B6' = 'География - Room 109 - Teacher: Марина Г.Б.'
    'C6' = 'Математика - Room 205 - Teacher: Константин В.Л.'
    'D6' = 'Химия - Room 212 - Teacher: Оксана М.С.'
    'E6' = 'Литература - Room 212 - Teacher: Георгий Н.М.'
    'F6' = 'Физика - Room 110 - Teacher: Роман А.Д.'
    'G6' = 'Информатика - Room 107 - Teacher: Анастасия А.М.'
    'H6' = 'Информатика - Room 107 - Teacher: Анастасия А.М.'
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet '11A' ----
$ws = $wb.Worksheets.Item('11A')

# Add Lesson 6 / Lesson 7 header cells, copying the header style from F1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G1").Value = "Lesson 6"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H1").Value = "Lesson 7"

$values = @{
    'B2' = 'Математика - Room 101 - Teacher: Евгения О.П.'
    'C2' = 'Математика - Room 205 - Teacher: Алексей С.К.'
    'D2' = 'История - Room 114 - Teacher: Елена В.К.'
    'E2' = 'Обществознание - Room 213 - Teacher: Наталья Д.З.'
    'F2' = 'Математика - Room 212 - Teacher: Марина Г.Б.'
    'G2' = 'Физическая культура - Room 111 - Teacher: Светлана А.К.'
    'H2' = 'Русский язык - Room 101 - Teacher: Михаил С.К.'
    'B3' = 'Физическая культура - Room 111 - Teacher: Светлана А.К.'
    'C3' = 'Русский язык - Room 101 - Teacher: Екатерина А.П.'
    'D3' = 'X - Room 0 - Teacher: X'
    'E3' = 'Иностранный язык - Room 206 - Teacher: Юлия Н.В.'
    'F3' = 'Обществознание - Room 212 - Teacher: Наталья Д.З.'
    'G3' = 'Физика - Room 110 - Teacher: Роман А.Д.'
    'H3' = 'Иностранный язык - Room 206 - Teacher: Елена В.К.'
    'B4' = 'Литература - Room 202 - Teacher: Илья В.М.'
    'C4' = 'Физика - Room 110 - Teacher: Игорь В.Н.'
    'D4' = 'Химия - Room 211 - Teacher: Елена С.Т.'
    'E4' = 'География - Room 109 - Teacher: Марина Г.Б.'
    'F4' = 'Математика - Room 109 - Teacher: Алексей С.К.'
    'G4' = 'Биология - Room 310 - Teacher: Андрей П.С.'
    'H4' = 'Химия - Room 212 - Teacher: Оксана М.С.'
    'B5' = 'Информатика - Room 107 - Teacher: Алёна В.Л.'
    'C5' = 'Литература - Room 202 - Teacher: Владимир В.П.'
    'D5' = 'История - Room 111 - Teacher: Елена В.К.'
    'E5' = 'Химия - Room 209 - Teacher: Елена С.Т.'
    'F5' = 'Литература - Room 207 - Teacher: Владимир В.П.'
    'G5' = 'Информатика - Room 107 - Teacher: Алёна В.Л.'
    'H5' = 'География - Room 111 - Teacher: Василий И.М.'
    'B6' = 'Биология - Room 310 - Teacher: Светлана А.К.'
    'C6' = 'Физика - Room 111 - Teacher: Игорь В.Н.'
    'D6' = 'Русский язык - Room 213 - Teacher: Александр И.Р.'
    'E6' = 'Математика - Room 109 - Teacher: Сергей А.Т.'
    'F6' = 'Русский язык - Room 110 - Teacher: Анна С.И.'
    'G6' = 'Математика - Room 109 - Teacher: Константин В.Л.'
    'H6' = 'Информатика - Room 107 - Teacher: Алёна В.Л.'
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet '11B' ----
$ws = $wb.Worksheets.Item('11B')

# Add Lesson 6 / Lesson 7 header cells, copying the header style from F1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G1").Value = "Lesson 6"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H1").Value = "Lesson 7"

$values = @{
    'B2' = 'Обществознание - Room 214 - Teacher: Анна С.И.'
    'C2' = 'Литература - Room 210 - Teacher: Владимир В.П.'
    'D2' = 'Информатика - Room 107 - Teacher: Игорь В.Н.'
    'E2' = 'Русский язык - Room 101 - Teacher: Дмитрий А.Ф.'
    'F2' = 'Литература - Room 202 - Teacher: Владимир В.П.'
    'G2' = 'Иностранный язык - Room 206 - Teacher: Елена В.К.'
    'H2' = 'Информатика - Room 107 - Teacher: Игорь В.Н.'
    'B3' = 'Математика - Room 104 - Teacher: Алексей С.К.'
    'C3' = 'Математика - Room 212 - Teacher: Евгения О.П.'
    'D3' = 'Иностранный язык - Room 208 - Teacher: Елена В.К.'
    'E3' = 'Математика - Room 109 - Teacher: Марина Г.Б.'
    'F3' = 'Химия - Room 213 - Teacher: Оксана М.С.'
    'G3' = 'История - Room 111 - Teacher: Наталья Д.З.'
    'H3' = 'Математика - Room 109 - Teacher: Евгения О.П.'
    'B4' = 'Биология - Room 310 - Teacher: Светлана А.К.'
    'C4' = 'Математика - Room 109 - Teacher: Евгения О.П.'
    'D4' = 'Русский язык - Room 206 - Teacher: Михаил С.К.'
    'E4' = 'Русский язык - Room 213 - Teacher: Дмитрий А.Ф.'
    'F4' = 'География - Room 111 - Teacher: Марина Г.Б.'
    'G4' = 'Физика - Room 114 - Teacher: Роман А.Д.'
    'H4' = 'Математика - Room 205 - Teacher: Евгения О.П.'
    'B5' = 'X - Room 0 - Teacher: X'
    'C5' = 'Обществознание - Room 214 - Teacher: Анна С.И.'
    'D5' = 'Физическая культура - Room 112 - Teacher: Светлана А.К.'
    'E5' = 'Физическая культура - Room 112 - Teacher: Юлия Н.В.'
    'F5' = 'Литература - Room 202 - Teacher: Илья В.М.'
    'G5' = 'Химия - Room 211 - Teacher: Оксана М.С.'
    'H5' = 'Информатика - Room 107 - Teacher: Анастасия А.М.'
    'B6' = 'Физика - Room 114 - Teacher: Роман А.Д.'
    'C6' = 'Русский язык - Room 206 - Teacher: Михаил С.К.'
    'D6' = 'Физика - Room 110 - Teacher: Роман А.Д.'
    'E6' = 'География - Room 111 - Teacher: Марина Г.Б.'
    'F6' = 'Биология - Room 309 - Teacher: Светлана А.К.'
    'G6' = 'Химия - Room 213 - Teacher: Оксана М.С.'
    'H6' = 'История - Room 114 - Teacher: Наталья Д.З.'
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet '11C' ----
$ws = $wb.Worksheets.Item('11C')

# Add Lesson 6 / Lesson 7 header cells, copying the header style from F1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G1").Value = "Lesson 6"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H1").Value = "Lesson 7"

$values = @{
    'B2' = 'Русский язык - Room 213 - Teacher: Александр И.Р.'
    'C2' = 'Физическая культура - Room 109 - Teacher: Светлана А.К.'
    'D2' = 'Русский язык - Room 213 - Teacher: Ольга В.С.'
    'E2' = 'Литература - Room 210 - Teacher: Илья В.М.'
    'F2' = 'Русский язык - Room 204 - Teacher: Екатерина А.П.'
    'G2' = 'Русский язык - Room 108 - Teacher: Михаил С.К.'
    'H2' = 'Иностранный язык - Room 208 - Teacher: Елена В.К.'
    'B3' = 'География - Room 111 - Teacher: Марина Г.Б.'
    'C3' = 'Математика - Room 108 - Teacher: Марина Г.Б.'
    'D3' = 'Информатика - Room 107 - Teacher: Алёна В.Л.'
    'E3' = 'Математика - Room 109 - Teacher: Евгения О.П.'
    'F3' = 'Физика - Room 113 - Teacher: Игорь В.Н.'
    'G3' = 'Физика - Room 111 - Teacher: Роман А.Д.'
    'H3' = 'Обществознание - Room 212 - Teacher: Наталья Д.З.'
    'B4' = 'Физика - Room 113 - Teacher: Игорь В.Н.'
    'C4' = 'Математика - Room 108 - Teacher: Константин В.Л.'
    'D4' = 'X - Room 0 - Teacher: X'
    'E4' = 'Литература - Room 103 - Teacher: Петр И.С.'
    'F4' = 'География - Room 109 - Teacher: Василий И.М.'
    'G4' = 'Физическая культура - Room 109 - Teacher: Светлана А.К.'
    'H4' = 'Биология - Room 310 - Teacher: Андрей П.С.'
    'B5' = 'Химия - Room 212 - Teacher: Анастасия А.Ф.'
    'C5' = 'Химия - Room 212 - Teacher: Оксана М.С.'
    'D5' = 'История - Room 110 - Teacher: Наталья Д.З.'
    'E5' = 'Литература - Room 105 - Teacher: Петр И.С.'
    'F5' = 'Иностранный язык - Room 205 - Teacher: Юлия Н.В.'
    'G5' = 'История - Room 110 - Teacher: Елена В.К.'
    'H5' = 'Информатика - Room 107 - Teacher: Алёна В.Л.'
    'B6' = 'Математика - Room 104 - Teacher: Евгения О.П.'
    'C6' = 'Математика - Room 108 - Teacher: Евгения О.П.'
    'D6' = 'Информатика - Room 107 - Teacher: Игорь В.Н.'
    'E6' = 'Химия - Room 213 - Teacher: Оксана М.С.'
    'F6' = 'Обществознание - Room 213 - Teacher: Наталья Д.З.'
    'G6' = 'Математика - Room 109 - Teacher: Константин В.Л.'
    'H6' = 'Биология - Room 310 - Teacher: Андрей П.С.'
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

